# Update the cash-flow data row (row 2) on the active sheet to reflect the
# 2017-12-31 report period figures instead of the 2019-12-31 ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report date
$ws.Range("N2").Value = "2017-12-31 00:00:00"

# Net cash flow metrics + their YoY ratios
$ws.Range("O2").Value = -10381535.22
$ws.Range("P2").Value = -82.68743924419999
$ws.Range("Q2").Value = 185328536.36
$ws.Range("R2").Value = 1476.1152147282
$ws.Range("S2").Value = 66961358.53
$ws.Range("T2").Value = 533.3376179748
$ws.Range("U2").Value = -30135332
$ws.Range("V2").Value = -240.0235977674

# RECEIVE_INVEST_INCOME / RII_RATIO have no reported value for this period
$ws.Range("W2").ClearContents()
$ws.Range("X2").ClearContents()

$ws.Range("Y2").Value = 11386621.6
$ws.Range("Z2").Value = 90.69280812460001
$ws.Range("AA2").Value = 28206889.16
$ws.Range("AB2").Value = 224.6638270987
$ws.Range("AC2").Value = -12555153.86

# CCE_ADD_RATIO has no reported value for this period
$ws.Range("AD2").ClearContents()
